$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'326.65"
$ws.Range("D3").Value = "'44.17"
$ws.Range("E3").Value = "'0.41%"
$ws.Range("D4").Value = "'5.565"
$ws.Range("E4").Value = "'-3.11%"
$ws.Range("D5").Value = "'0.08031"
$ws.Range("E5").Value = "'-3.99%"
$ws.Range("D6").Value = "'4.297"
$ws.Range("E6").Value = "'-5.00%"
$ws.Range("D7").Value = "'1.895"
$ws.Range("E7").Value = "'-2.86%"
$ws.Range("E8").Value = "'-8.93%"
$ws.Range("D9").Value = "'0.9446"
$ws.Range("E9").Value = "'-0.33%"
$ws.Range("D10").Value = "'0.1155"
$ws.Range("E10").Value = "'-6.95%"
$ws.Range("E11").Value = "'-6.99%"
$ws.Range("D12").Value = "'0.09694"
$ws.Range("E12").Value = "'-2.44%"
$ws.Range("D13").Value = "'0.04377"
$ws.Range("E13").Value = "'-0.91%"
$ws.Range("D14").Value = "'0.1065"
$ws.Range("E14").Value = "'-0.26%"
$ws.Range("D15").Value = "'0.001287"
$ws.Range("E15").Value = "'-0.60%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005982"
$ws.Range("E16").Value = "'-1.40%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.622"
$ws.Range("E17").Value = "'3.95%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "'0.3495"
$ws.Range("E18").Value = "'-1.19%"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "'8.596"
$ws.Range("E19").Value = "'-1.26%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1390"
$ws.Range("E20").Value = "'1.91%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "'0.2529"
$ws.Range("E21").Value = "'-3.30%"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "'0.04221"
$ws.Range("E22").Value = "'-4.14%"
$ws.Range("E23").Value = "'0.17%"
$ws.Range("D24").Value = "'0.004480"
$ws.Range("E24").Value = "'2.71%"
$ws.Range("D25").Value = "'0.0001263"
$ws.Range("E25").Value = "'-0.05%"
$ws.Range("D26").Value = "'0.0003997"
$ws.Range("D38").Value = "'0.02613"
$ws.Range("E38").Value = "'-7.35%"
$ws.Range("D39").Value = "'0.05427"
$ws.Range("E39").Value = "'-7.91%"
$ws.Range("D40").Value = "'0.007598"
$ws.Range("E40").Value = "'-4.60%"
$ws.Range("E41").Value = "'-2.14%"
$ws.Range("D42").Value = "'0.007269"
$ws.Range("E42").Value = "'-19.61%"
$ws.Range("D44").Value = "'0.008832"
$ws.Range("E44").Value = "'-15.09%"
$ws.Range("D45").Value = "'0.00006920"
$ws.Range("E45").Value = "'-4.35%"
$ws.Range("D46").Value = "'0.00000000752"
$ws.Range("E46").Value = "'-0.05%"
$ws.Range("D47").Value = "'0.003558"
$ws.Range("E47").Value = "'11.28%"
$ws.Range("D48").Value = "'0.002275"
$ws.Range("E49").Value = "'-0.05%"
$ws.Range("E50").Value = "'-0.05%"
